$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.028.40"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("D3").Value = "1.650.11"
$ws.Range("E3").Value = "  -5.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.60"
$ws.Range("E5").Value = "  -5.59%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").Value = "  -6.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2617"
$ws.Range("E8").Value = "  -5.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06016"
$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07192"
$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("D11").Value = "1.647.41"
$ws.Range("E11").Value = "  -5.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("E12").Value = "  -2.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6223"
$ws.Range("E13").Value = "  -4.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.561"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.04"
$ws.Range("E15").Value = "  -6.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "25.013.67"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.53"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006628"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.561"
$ws.Range("E21").Value = "  +6.35%  "

$ws.Range("D22").Value = "1.860.76"
$ws.Range("E22").Value = "  -5.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.627"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.316"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "132.22"
$ws.Range("E25").Value = "  -2.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.91"
$ws.Range("E26").Value = "  -2.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.398"
$ws.Range("E27").Value = "  -7.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.07"
$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.679"
$ws.Range("E29").Value = "  -5.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.766"
$ws.Range("E30").Value = "  -5.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07892"
$ws.Range("E31").Value = "  -4.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.591"
$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04545"
$ws.Range("E33").Value = "  -3.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9991"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.590"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9356"
$ws.Range("E36").Value = "  -6.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5801"
$ws.Range("E37").Value = "  -7.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.588"
$ws.Range("E38").Value = "  -5.26%  "

$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8440"
$ws.Range("E40").Value = "  +10.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9997"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.823"
$ws.Range("E42").Value = "  -5.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.54"
$ws.Range("E43").Value = "  -2.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3740"
$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.792"
$ws.Range("E45").Value = "  -4.71%  "

$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.142"
$ws.Range("E47").Value = "  -3.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05199"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.88"
$ws.Range("E49").Value = "  -3.01%  "

$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.51"
$ws.Range("E51").Value = "  -9.24%  "
